# Updated: po 20. 12. 2021
# Applies corrected AgTests (F) / AgPosit (G) figures for rows 628-652 and
# appends three new daily rows (653-655) with data through 2021-12-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (column F = AgTests, column G = AgPosit) ---
$ws.Cells.Item(628, 6).Value = 64272

$ws.Cells.Item(632, 6).Value = 44049

$ws.Cells.Item(635, 6).Value = 82773

$ws.Cells.Item(637, 6).Value = 43319

$ws.Cells.Item(638, 6).Value = 37324
$ws.Cells.Item(638, 7).Value = 1950

$ws.Cells.Item(639, 6).Value = 39870
$ws.Cells.Item(639, 7).Value = 1948

$ws.Cells.Item(640, 6).Value = 19446
$ws.Cells.Item(640, 7).Value = 1211

$ws.Cells.Item(641, 6).Value = 33528
$ws.Cells.Item(641, 7).Value = 1367

$ws.Cells.Item(642, 6).Value = 66955
$ws.Cells.Item(642, 7).Value = 2372

$ws.Cells.Item(643, 6).Value = 42597
$ws.Cells.Item(643, 7).Value = 1639

$ws.Cells.Item(644, 6).Value = 36149
$ws.Cells.Item(644, 7).Value = 1471

$ws.Cells.Item(645, 6).Value = 35118
$ws.Cells.Item(645, 7).Value = 1288

$ws.Cells.Item(646, 6).Value = 35659
$ws.Cells.Item(646, 7).Value = 1337

$ws.Cells.Item(647, 6).Value = 16041
$ws.Cells.Item(647, 7).Value = 898

$ws.Cells.Item(648, 6).Value = 28799
$ws.Cells.Item(648, 7).Value = 1029

$ws.Cells.Item(649, 6).Value = 61429
$ws.Cells.Item(649, 7).Value = 1787

$ws.Cells.Item(650, 6).Value = 36915
$ws.Cells.Item(650, 7).Value = 1162

$ws.Cells.Item(651, 6).Value = 34558
$ws.Cells.Item(651, 7).Value = 1024

$ws.Cells.Item(652, 6).Value = 33722
$ws.Cells.Item(652, 7).Value = 1019

# --- New rows appended at the bottom (2021-12-17 through 2021-12-19) ---
$newRows = @(
    @{ Row = 653; A = 44547; B = 802684; C = 14252; D = 4021; E = 15931; F = 29246; G = 836 },
    @{ Row = 654; A = 44548; B = 806115; C = 11383; D = 3431; E = 16014; F = 10984; G = 514 },
    @{ Row = 655; A = 44549; B = 807339; C = 4713;  D = 1224; E = 16068; F = 15316; G = 475 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
